# Update the "Periodo Mora" value from 2507 to 2508 for all worker rows
# (base de datos EC actualizada a nuevo periodo de mora).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16:E20").Value = "2508"
